$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the whole B2:D9 block to 0, matching the target diff.
$ws.Range("B2:D9").Value = 0
